$d = $word.ActiveDocument

# The first two paragraphs of the document are:
#   1) "On Pilgrimage - June 1953"  (Heading1)
#   2) "By Dorothy Day"             (bold run)
# They need to become:
#   1) "On Pilgrimage - June 1953"  (Title style, split across several runs)
#   2) "Dorothy Day"                (Authors style, split across several runs)
# wrapped by a bookmark that spans paragraph 1 only (the bookmark itself is
# left untouched by this script - the COM bridge does not expose bookmark
# mutation - but its surrounding content is rewritten in place).

$titlePara = $d.Paragraphs.Item(1)
$authorPara = $d.Paragraphs.Item(2)

$startPos = $titlePara.Range.Start
$endPos = $authorPara.Range.End

$targetRange = $d.Range($startPos, $endPos)

$newXml = @'
<?xml version="1.0" encoding="utf-8"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Title"/>
            </w:pPr>
            <w:r><w:t xml:space="preserve">On</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:t xml:space="preserve">Pilgrimage</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:t xml:space="preserve">-</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:t xml:space="preserve">June</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:t xml:space="preserve">1953</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Authors"/>
            </w:pPr>
            <w:r><w:t xml:space="preserve">Dorothy</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:t xml:space="preserve">Day</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$targetRange.InsertXML($newXml)
